{"js": "// Split the \"Write the SQL that creates the tables, indexes and sequences for\n// your database...\" sentence into several runs: highlight \"sequences\" in\n// yellow and insert \"(for userID)\" right after it (with the matching\n// proofErr grammar/spelling markers Word leaves behind), while leaving the\n// rest of the sentence's text/formatting untouched.\n\nconst originalSentence =\n  \"Write the SQL that creates the tables, indexes and sequences for your \" +\n  \"database. (Place all your SQL in the same script file.)\";\n\nconst searchResults = context.document.body.search(originalSentence, {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target sentence in the document.\");\n}\n\nconst targetRange = searchResults.items[0];\n\n// Flat-OPC wrapped OOXML for the replacement runs (required shape for\n// Range.insertOoxml / Body.insertOoxml). Replacing the whole matched range\n// keeps the paragraph (w:pPr) intact and swaps only its run content.\nconst replacementOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidRPr=\"002027F2\">\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Write the SQL that creates the tables, indexes and </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:highlight w:val=\"yellow\"/>\n              </w:rPr>\n              <w:t>sequences</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t>(</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">for </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t>userID</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t>)</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> for your database. (Place all your SQL in the same script file.)</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Split the \"Write the SQL that creates the tables, indexes and sequences for\n# your database...\" sentence into several runs: highlight \"sequences\" in\n# yellow and insert \"(for userID)\" right after it (with the matching\n# proofErr grammar/spelling markers Word leaves behind), while leaving the\n# rest of the sentence's text/formatting untouched.\n\n$d = $word.ActiveDocument\n\n$originalSentence = \"Write the SQL that creates the tables, indexes and sequences for your database. (Place all your SQL in the same script file.)\"\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($originalSentence)\nif (-not $found) {\n    throw \"Could not find the target sentence in the document.\"\n}\n\n# Re-seat a plain Range over the matched span so InsertXML REPLACES that\n# range's contents instead of inserting next to a Find-owned range.\n$targetRange = $d.Range($searchRange.Start, $searchRange.End)\n\n$replacementOoxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidRPr=\"002027F2\">\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Write the SQL that creates the tables, indexes and </w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:highlight w:val=\"yellow\"/>\n              </w:rPr>\n              <w:t>sequences</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t>(</w:t>\n            </w:r>\n            <w:proofErr w:type=\"gramEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">for </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t>userID</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t>)</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> for your database. (Place all your SQL in the same script file.)</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$targetRange.InsertXML($replacementOoxml)\n"}
